$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the empty "grandes regiões e unidades da federação" row (row 6).
# This shifts all subsequent rows up by one and removes the now-unused
# shared string from the string table, and removes the last (now blank)
# row that was pushed off the bottom (former row 37).
$ws.Rows.Item(6).Delete()
